$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "301.13"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.80%"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.50%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.996"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.88%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07689"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.41%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.083"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-8.04%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.914"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-2.06%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.031"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.48%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9141"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.52%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09670"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "7.95%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1865"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.65%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08519"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.60%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03538"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-6.41%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09961"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.23%"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.48%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005675"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.28%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.59%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "11.47%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1329"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.52%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.756"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "4.29%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-1.62%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04584"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-2.00%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.005093"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "12.46%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001232"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.16%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001402"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "7.54%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01756"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-0.62%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04618"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-2.41%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007486"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-6.18%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-1.60%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007730"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.32%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002242"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-2.74%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01028"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "6.82%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006290"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.42%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.13%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0005802"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "0.03%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "38.05"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "555.71%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.002002"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-25.75%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.00002103"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.13%"
